# Remove the leftover "Google Shape;88;p16" (date placeholder, ph type="dt")
# and "Google Shape;89;p16" (footer placeholder, ph type="ftr") text boxes
# from the last three slides (11, 12, 13). These were stray per-slide
# date/"Consultora: COVENAR" footer shapes left over from the Google Slides
# import; the author removed them in this commit.
#
# NOTE: calling .Delete() on these shapes only clears/resets the
# placeholder (PowerPoint keeps required placeholder shapes around), so we
# use .Cut() instead, which actually removes the shape node from the slide.

$p = $ppt.ActivePresentation

foreach ($slideIdx in 11, 12, 13) {
    $s = $p.Slides.Item($slideIdx)

    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shape = $s.Shapes.Item($i)
        $name = $shape.Name
        if ($name -eq "Google Shape;88;p16" -or $name -eq "Google Shape;89;p16") {
            $shape.Cut()
        }
    }
}
